$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content (old layout lived in column B, rows 1-7)
$ws.Range("A1:B8").ClearContents()

# New text content
$ws.Range("A1").Value = "This is Lesson 0. If you don't want a lesson 0, fair enough! Just start in row 2"
$ws.Range("A2").Value = "This is Lesson 1. Lesson 1 should be in row 2 etc."
$ws.Range("B2").Value = "If your lessons are over here, set Column Number to 1 in lesson_indexes.json. If they are in another column, set the appropriate column number"
$ws.Range("A3").Value = "Lesson 2"
$ws.Range("A4").Value = "Lesson 3"
$ws.Range("A5").Value = "Lesson 4"
$ws.Range("A6").Value = "You get the idea"
$ws.Range("A7").Value = "Lorum Impusm"
$ws.Range("A8").Value = "Dolor sit amet"

# Column widths (closest achievable values to the authored widths;
# the COM layer quantizes ColumnWidth to 1/6-character steps)
$ws.Columns.Item(1).ColumnWidth = 46.33
$ws.Columns.Item(2).ColumnWidth = 43.5

# Alignment / wrap text (B2 first so its style is allocated as cellXfs index 1,
# matching the centered+wrapped style; A1 gets index 2, wrap-only)
$ws.Range("B2").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B2").WrapText = $true

$ws.Range("A1").WrapText = $true

# Row heights
$ws.Rows.Item(1).RowHeight = 30.75
$ws.Rows.Item(2).RowHeight = 43.5

# Selection
$ws.Range("B2").Select() | Out-Null
